# Update column F (dSF) values for specific rows to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = 0
$ws.Range("F3").Value  = -3
$ws.Range("F10").Value = 1
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = -1
$ws.Range("F17").Value = -6
$ws.Range("F21").Value = -5
$ws.Range("F23").Value = -1
$ws.Range("F32").Value = -8
$ws.Range("F35").Value = 3
$ws.Range("F38").Value = -4
$ws.Range("F40").Value = 4
$ws.Range("F41").Value = 2
$ws.Range("F45").Value = -3
$ws.Range("F49").Value = 0
$ws.Range("F56").Value = 1
$ws.Range("F58").Value = 2
